$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A18:D39").ClearContents()
